# Append 19 new "한울소재과학" rows (rows 24-42) to the records sheet.
# Columns: A=회사명, B=발행시간, C=회차 (numeric), D=추가주식수(주) (text),
#          E=발행/전환/행사가액(원) (text). D/E keep thousands-separator
#          text formatting, so we force the cell to Text before writing the
#          value (otherwise Excel auto-parses "217,013" as the number
#          217013) and then reset the style back to Normal so no stray
#          NumberFormat-derived style index is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("한울소재과학", "2025-10-02 17:05", 4, "217,013",   "2,304"),
    @("한울소재과학", "2025-10-02 17:05", 5, "217,013",   "2,304"),
    @("한울소재과학", "2025-07-09 16:45", 5, "173,611",   "2,304"),
    @("한울소재과학", "2025-07-09 16:45", 5, "108,506",   "2,304"),
    @("한울소재과학", "2025-07-01 17:40", 2, "217,013",   "2,304"),
    @("한울소재과학", "2025-07-01 17:39", 5, "520,832",   "2,304"),
    @("한울소재과학", "2025-05-28 16:39", 2, "86,805",    "2,304"),
    @("한울소재과학", "2025-04-23 16:58", 4, "651,039",   "2,304"),
    @("한울소재과학", "2025-04-23 16:58", 5, "325,519",   "2,304"),
    @("한울소재과학", "2025-04-16 16:46", 2, "781,245",   "2,304"),
    @("한울소재과학", "2025-04-16 16:39", 5, "260,416",   "2,304"),
    @("한울소재과학", "2025-03-26 17:11", 5, "108,506",   "2,304"),
    @("한울소재과학", "2025-03-20 17:23", 4, "651,040",   "2,304"),
    @("한울소재과학", "2025-03-20 17:23", 5, "173,610",   "2,304"),
    @("한울소재과학", "2025-03-14 17:24", 4, "217,012",   "2,304"),
    @("한울소재과학", "2025-03-14 17:24", 5, "43,402",    "2,304"),
    @("한울소재과학", "2025-02-19 17:37", 4, "4,014,741", "2,304"),
    @("한울소재과학", "2025-02-19 17:37", 5, "1,866,309", "2,304"),
    @("한울소재과학", "2025-02-19 17:36", 3, "57,502",    "10,869")
)

$startRow = 24
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]

    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 5).Style = "Normal"
}
